$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 170; existing rows 170-207 shift down to 172-209.
$ws.Rows("170:171").Insert()

# --- New row 170 ---
$ws.Cells.Item(170, 1).Value = 4
$ws.Cells.Item(170, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(170, 3).Value = "Los Lagos"
$ws.Cells.Item(170, 4).Value = "3/22/2022"
$ws.Cells.Item(170, 5).Value = 10
$ws.Cells.Item(170, 6).Value = "Fruta"
$ws.Cells.Item(170, 7).Value = 100103
$ws.Cells.Item(170, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(170, 9).Value = 100103004
$ws.Cells.Item(170, 10).Value = "Durazno"
$ws.Cells.Item(170, 11).Value = "September Snow"
$ws.Cells.Item(170, 12).Value = "Especial"
$ws.Cells.Item(170, 13).Value = 100
$ws.Cells.Item(170, 14).Value = 20000
$ws.Cells.Item(170, 15).Value = 20000
$ws.Cells.Item(170, 16).Value = 20000
$ws.Cells.Item(170, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(170, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(170, 19).Value = 1333
$ws.Cells.Item(170, 20).Value = 15

# --- New row 171 ---
$ws.Cells.Item(171, 1).Value = 4
$ws.Cells.Item(171, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(171, 3).Value = "Los Lagos"
$ws.Cells.Item(171, 4).Value = "3/22/2022"
$ws.Cells.Item(171, 5).Value = 10
$ws.Cells.Item(171, 6).Value = "Fruta"
$ws.Cells.Item(171, 7).Value = 100103
$ws.Cells.Item(171, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(171, 9).Value = 100103004
$ws.Cells.Item(171, 10).Value = "Durazno"
$ws.Cells.Item(171, 11).Value = "September Snow"
$ws.Cells.Item(171, 12).Value = "Primera"
$ws.Cells.Item(171, 13).Value = 200
$ws.Cells.Item(171, 14).Value = 16000
$ws.Cells.Item(171, 15).Value = 17000
$ws.Cells.Item(171, 16).Value = 16500
$ws.Cells.Item(171, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(171, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(171, 19).Value = 1100
$ws.Cells.Item(171, 20).Value = 15
